$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Append the three new rows (construct-mapping entries) at the bottom ---
$newRows = @(
    @{
        B = "mass media mode of delivery"
        C = "Informational mode of delivery of radio broadcast, television, online press and printed press to a mass audience."
        F = "informational mode of delivery"
        L = "This is a form of mass media."
    },
    @{
        B = "online press mode of delivery"
        C = "Electronic mode of delivery of a newspaper or magazine."
        F = "electronic mode of delivery"
    },
    @{
        B = "social influence intervention through mass media"
        C = "A behaviour change intervention that is an awareness of other people’s thoughts, feelings and actions BCT delivered through a mass media mode of delivery."
        F = "behaviour change intervention"
    }
)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$templateRow = $ws.Range("A2:V2")
$scratch = $ws.Range("Z1")

foreach ($row in $newRows) {
    $lastRow = $lastRow + 1
    $destRow = $ws.Range("A" + $lastRow + ":V" + $lastRow)

    # Copy row 2's formatting (style index), then its values (the "Process" /
    # "Mode of delivery" / "Published" / "0" constants shared by every data row).
    $templateRow.Copy() | Out-Null
    $destRow.PasteSpecial(-4122) | Out-Null
    $ws.Application.CutCopyMode = $false

    $templateRow.Copy() | Out-Null
    $destRow.PasteSpecial(-4163) | Out-Null
    $ws.Application.CutCopyMode = $false

    # Columns A, D, E are blank for the new rows.
    $ws.Cells.Item($lastRow, 1).ClearContents() | Out-Null
    $ws.Cells.Item($lastRow, 4).ClearContents() | Out-Null
    $ws.Cells.Item($lastRow, 5).ClearContents() | Out-Null

    # Write the differing text cells through a scratch cell + paste-special (values
    # only) so the destination keeps its style (direct .Value assignment here would
    # otherwise reset the cell's style index in this runtime).
    $scratch.Value = $row.B
    $scratch.Copy() | Out-Null
    $ws.Cells.Item($lastRow, 2).PasteSpecial(-4163) | Out-Null
    $ws.Application.CutCopyMode = $false

    $scratch.Value = $row.C
    $scratch.Copy() | Out-Null
    $ws.Cells.Item($lastRow, 3).PasteSpecial(-4163) | Out-Null
    $ws.Application.CutCopyMode = $false

    $scratch.Value = $row.F
    $scratch.Copy() | Out-Null
    $ws.Cells.Item($lastRow, 6).PasteSpecial(-4163) | Out-Null
    $ws.Application.CutCopyMode = $false

    if ($row.ContainsKey("L")) {
        $scratch.Value = $row.L
        $scratch.Copy() | Out-Null
        $ws.Cells.Item($lastRow, 12).PasteSpecial(-4163) | Out-Null
        $ws.Application.CutCopyMode = $false
    } else {
        $ws.Cells.Item($lastRow, 12).ClearContents() | Out-Null
    }

    $scratch.ClearContents() | Out-Null
}

# --- 2. Sort the whole data range alphabetically by column B (Label) ---
$dataRange = $ws.Range("A2:V" + $lastRow)
$sortKey = $ws.Range("B2:B" + $lastRow)
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($sortKey, 0, 1, 0, 0) | Out-Null
$ws.Sort.SetRange($dataRange)
$ws.Sort.Header = 0
$ws.Sort.Apply()

# --- 3. AutoFit column B width ---
$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null

# --- 4. Update the active view state (scroll + selection) ---
$ws.Application.ActiveWindow.ScrollRow = 39
$ws.Range("V58").Select() | Out-Null
